# Bring back the "Tipo"/"valor" input column on both sheets (Maquina + Humano),
# and make "Maquina" the active sheet/tab again.

$wb = $excel.ActiveWorkbook

# --- Sheet "Maquina": add column F ("Tipo" header, "valor" for every data row) ---
$wsM = $wb.Worksheets.Item("Maquina")
$wsM.Activate()

$wsM.Range("F1").Value = "Tipo"
for ($r = 2; $r -le 32; $r++) {
    $wsM.Cells.Item($r, 6).Value = "valor"
}

$wsM.Range("F1").Select()
$wsM.Range("F1:F32").Select()

# --- Sheet "Humano": add column E ("Tipo" header, "valor" for every data row) ---
$wsH = $wb.Worksheets.Item("Humano")

$wsH.Range("E1").Value = "Tipo"
for ($r = 2; $r -le 32; $r++) {
    $wsH.Cells.Item($r, 5).Value = "valor"
}

$wsH.Range("E1").Select()

# Re-activate Maquina so it is the active/selected sheet at save time
$wsM.Activate()
